$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was added for "Papa" at Feria Lagunitas de
# Puerto Montt. It lands at row 396, pushing the existing rows 396-415 down
# to 397-416 (dimension grows from R415 to R416).
$ws.Rows(396).Insert()

$ws.Range("A396").Value = 4
$ws.Range("B396").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C396").Value = "Los Lagos"
$ws.Range("D396").Value = 44753
$ws.Range("E396").Value = 10
$ws.Range("F396").Value = 100114001
$ws.Range("G396").Value = "Papa"
$ws.Range("H396").Value = "Patagonia"
$ws.Range("I396").Value = "1a (guarda)"
$ws.Range("J396").Value = 250
$ws.Range("K396").Value = 8000
$ws.Range("L396").Value = 8000
$ws.Range("M396").Value = 8000
$ws.Range("N396").Value = '$/saco 25 kilos'
$ws.Range("O396").Value = "Provincia de Llanquihue"
$ws.Range("P396").Value = 320
$ws.Range("Q396").Value = 25
$ws.Range("R396").Value = "Hortaliza"
